# fitur jika ada score yang sama
# Reorder tied rows (same DeFuzzy Score) and update one row whose
# underlying fuzzy inputs/score changed.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows 2-5 all share the same DeFuzzy Score (83.33333333333333) and get
# re-ranked among themselves when scores tie.
$ws.Range("A2").Value = 24
$ws.Range("B2").Value = 100
$ws.Range("C2").Value = 9

$ws.Range("A3").Value = 42
$ws.Range("B3").Value = 94
$ws.Range("C3").Value = 10

$ws.Range("A4").Value = 79
$ws.Range("B4").Value = 87
$ws.Range("C4").Value = 9

$ws.Range("A5").Value = 69
$ws.Range("B5").Value = 86
$ws.Range("C5").Value = 10

# Row 11's underlying values/score change (new data / tie-break recalculation).
$ws.Range("A11").Value = 54
$ws.Range("B11").Value = 64
$ws.Range("C11").Value = 10
$ws.Range("D11").Value = 76.02689486552569
